$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the " QUALITY" row (row 3). This shifts MATERIAL, MATTERMATTER,
# ORGANSIM, PROTEIN, and TREE each up by one row, and the former last
# row (TREE, row 8) disappears, matching the diff exactly since the
# normalized-count and rank values are identical across rows 3-8.
$ws.Rows.Item(3).Delete()
